$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.968.17"
$ws.Range("E2").Value = "  -2.24%  "

$ws.Range("D3").Value = "1.649.18"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").Value = "'310.06"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "'0.3903"
$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("D8").Value = "'0.3809"
$ws.Range("E8").Value = "  -2.84%  "

$ws.Range("D9").Value = "'52.18"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").Value = "'1.348"
$ws.Range("E10").Value = "  -4.42%  "

$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "'0.08450"
$ws.Range("E12").Value = "  -1.90%  "

$ws.Range("D13").Value = "'23.86"
$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("D14").Value = "'7.060"
$ws.Range("E14").Value = "  -4.26%  "

$ws.Range("D15").Value = "'8.000"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("E16").Value = "  -2.87%  "

$ws.Range("D17").Value = "1.650.95"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").Value = "'94.24"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("D19").Value = "'0.06998"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").Value = "'19.68"
$ws.Range("E20").Value = "  -4.90%  "

$ws.Range("D21").Value = "'6.973"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "'13.76"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Value = "23.960.23"
$ws.Range("E24").Value = "  -2.19%  "

$ws.Range("D25").Value = "'2.445"
$ws.Range("E25").Value = "  +0.85%  "

$ws.Range("D26").Value = "'2.946"
$ws.Range("E26").Value = "  -3.29%  "

$ws.Range("D27").Value = "'22.05"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("D28").Value = "'152.86"
$ws.Range("E28").Value = "  -3.13%  "

$ws.Range("D29").Value = "'5.400"
$ws.Range("E29").Value = "  -1.76%  "

$ws.Range("D30").Value = "'137.88"
$ws.Range("E30").Value = "  -3.76%  "

$ws.Range("D31").Value = "'7.905"
$ws.Range("E31").Value = "  -3.10%  "

$ws.Range("D32").Value = "'2.514"
$ws.Range("E32").Value = "  -1.37%  "

$ws.Range("D33").Value = "1.829.65"
$ws.Range("E33").Value = "  -1.08%  "

$ws.Range("D34").Value = "'1.022"
$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("D35").Value = "'0.08045"
$ws.Range("E35").Value = "  -3.05%  "

$ws.Range("D36").Value = "'6.771"
$ws.Range("E36").Value = "  -1.61%  "

$ws.Range("D37").Value = "'0.02924"
$ws.Range("E37").Value = "  -3.87%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2676"
$ws.Range("E38").Value = "  -3.51%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'10.73"
$ws.Range("E39").Value = "  -4.74%  "

$ws.Range("D40").Value = "'0.09093"
$ws.Range("E40").Value = "  -1.85%  "

$ws.Range("D41").Value = "'0.7599"
$ws.Range("E41").Value = "  -2.45%  "

$ws.Range("D42").Value = "'13.38"
$ws.Range("E42").Value = "  -3.83%  "

$ws.Range("E43").Value = "  -1.83%  "

$ws.Range("D44").Value = "'16.31"
$ws.Range("E44").Value = "  -1.64%  "

$ws.Range("D45").Value = "'0.6973"
$ws.Range("E45").Value = "  -2.62%  "

$ws.Range("D46").Value = "'2.452"
$ws.Range("E46").Value = "  -3.83%  "

$ws.Range("D47").Value = "'4.088"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").Value = "'0.08319"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("D50").Value = "'134.87"
$ws.Range("E50").Value = "  -1.47%  "

$ws.Range("D51").Value = "'1.231"
$ws.Range("E51").Value = "  -4.05%  "
